$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the edited cells to Text format so numeric-looking / percent-looking
# strings are written back verbatim (matching the source inlineStr cells)
# instead of being auto-converted to numbers by Excel.
$editRange = $ws.Range("D2:E47")
$editRange.NumberFormat = "@"

$ws.Range("D2").Value = "287.44"
$ws.Range("E2").Value = "0.78%"
$ws.Range("E3").Value = "2.05%"
$ws.Range("D4").Value = "5.204"
$ws.Range("E4").Value = "2.04%"
$ws.Range("D5").Value = "0.06970"
$ws.Range("E5").Value = "4.40%"
$ws.Range("D6").Value = "7.423"
$ws.Range("E6").Value = "1.29%"
$ws.Range("D7").Value = "3.555"
$ws.Range("E7").Value = "4.95%"
$ws.Range("D8").Value = "1.404"
$ws.Range("E8").Value = "3.41%"
$ws.Range("D9").Value = "0.9002"
$ws.Range("E9").Value = "-3.83%"
$ws.Range("D10").Value = "0.1606"
$ws.Range("E10").Value = "2.04%"
$ws.Range("D11").Value = "0.07553"
$ws.Range("E11").Value = "16.59%"
$ws.Range("D12").Value = "0.07696"
$ws.Range("E12").Value = "0.72%"
$ws.Range("D13").Value = "0.02900"
$ws.Range("E13").Value = "0.93%"
$ws.Range("E14").Value = "0.49%"
$ws.Range("D15").Value = "0.001583"
$ws.Range("E15").Value = "-0.35%"
$ws.Range("D16").Value = "0.0006532"
$ws.Range("E16").Value = "1.71%"
$ws.Range("D17").Value = "0.006217"
$ws.Range("E17").Value = "1.52%"
$ws.Range("D18").Value = "3.489"
$ws.Range("E18").Value = "0.21%"
$ws.Range("E19").Value = "0.10%"
$ws.Range("D20").Value = "0.3245"
$ws.Range("E20").Value = "1.34%"
$ws.Range("D21").Value = "0.1337"
$ws.Range("E21").Value = "2.42%"
$ws.Range("D22").Value = "4.049"
$ws.Range("E22").Value = "-0.18%"
$ws.Range("D24").Value = "0.04531"
$ws.Range("E24").Value = "1.45%"
$ws.Range("E25").Value = "2.77%"
$ws.Range("E26").Value = "-7.36%"
$ws.Range("E27").Value = "-6.02%"
$ws.Range("E28").Value = "2.03%"
$ws.Range("D40").Value = "0.04378"
$ws.Range("E40").Value = "4.61%"
$ws.Range("D41").Value = "0.006941"
$ws.Range("E41").Value = "3.16%"
$ws.Range("D42").Value = "0.1247"
$ws.Range("E42").Value = "0.14%"
$ws.Range("D43").Value = "0.002070"
$ws.Range("E43").Value = "2.89%"
$ws.Range("D44").Value = "0.01177"
$ws.Range("E44").Value = "-2.31%"
$ws.Range("D45").Value = "0.00005841"
$ws.Range("E45").Value = "3.15%"
$ws.Range("E46").Value = "-1.85%"
$ws.Range("E47").Value = "-0.13%"

# Restore the default (unstyled) cell style so no stray formatting is introduced.
$editRange.Style = "Normal"

